{"js": "// Apply the three textual changes described by the diff:\n// 1) Consolidate the run-fragmented sentence \" This process should also\n//    include QA related activities such as fix verification, regression\n//    suit execution and continuous monitoring.\" into a single run (no\n//    wording change, just a clean re-write of the same visible text,\n//    mirroring the XML run-merge in the diff).\n// 2) Fix the typo \"pross platforms\" -> \"cross platforms\".\n// 3) Add a comma after \"In short\" -> \"In short, \".\n\nconst body = context.document.body;\n\n// --- Change 1: re-write the \"This process ... monitoring.\" sentence ---\nconst mergedSentence =\n  \" This process should also include QA related activities such as fix \" +\n  \"verification, regression suit execution and continuous monitoring.\";\nconst mergeResults = body.search(mergedSentence, { matchCase: true, matchWholeWord: false });\nmergeResults.load(\"items\");\nawait context.sync();\nfor (let i = 0; i < mergeResults.items.length; i++) {\n  mergeResults.items[i].insertText(mergedSentence, Word.InsertLocation.replace);\n}\nawait context.sync();\n\n// --- Change 2: \"pross platforms\" -> \"cross platforms\" ---\nconst typoResults = body.search(\"pross platforms\", { matchCase: true, matchWholeWord: false });\ntypoResults.load(\"items\");\nawait context.sync();\nfor (let i = 0; i < typoResults.items.length; i++) {\n  typoResults.items[i].insertText(\"cross platforms\", Word.InsertLocation.replace);\n}\nawait context.sync();\n\n// --- Change 3: \"In short \" -> \"In short, \" ---\nconst commaResults = body.search(\"In short \", { matchCase: true, matchWholeWord: false });\ncommaResults.load(\"items\");\nawait context.sync();\nfor (let i = 0; i < commaResults.items.length; i++) {\n  commaResults.items[i].insertText(\"In short, \", Word.InsertLocation.replace);\n}\nawait context.sync();\n", "ps1": "# Apply the three textual changes described by the diff:\n# 1) Re-write the run-fragmented sentence \" This process should also\n#    include QA related activities such as fix verification, regression\n#    suit execution and continuous monitoring.\" (no wording change, just\n#    a clean re-write of the same visible text that consolidates runs).\n# 2) Fix the typo \"pross platforms\" -> \"cross platforms\".\n# 3) Add a comma after \"In short\" -> \"In short, \".\n\n$d = $word.ActiveDocument\n\n# --- Change 1: re-write the \"This process ... monitoring.\" sentence ---\n$mergedSentence = \" This process should also include QA related activities such as fix verification, regression suit execution and continuous monitoring.\"\n$rng1 = $d.Range()\n$rng1.Find.Execute($mergedSentence, $false, $false, $false, $false, $false, $true, 1, $false, $mergedSentence, 2)\n\n# --- Change 2: \"pross platforms\" -> \"cross platforms\" ---\n$rng2 = $d.Range()\n$rng2.Find.Execute(\"pross platforms\", $false, $false, $false, $false, $false, $true, 1, $false, \"cross platforms\", 2)\n\n# --- Change 3: \"In short \" -> \"In short, \" ---\n$rng3 = $d.Range()\n$rng3.Find.Execute(\"In short \", $false, $false, $false, $false, $false, $true, 1, $false, \"In short, \", 2)\n"}
